$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 for the SMA connector itself (added first so shared strings are interned
# in the same order as the target: "SMA connector" then "CON-SMA-EDGE-S-ND")
$ws.Range("A15").Value = "SMA connector"
$ws.Range("B15").Value = "CON-SMA-EDGE-S-ND"
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 1.74
$ws.Range("E15").NumberFormat = $ws.Range("E14").NumberFormat

# Update the existing "Level shifter" part number (B14) to the new SMA connector's Digikey number
$ws.Range("B14").Value = "296-12163-1-ND"

# Column B auto-fits to the new, wider part number text (~18.78 chars in real Excel;
# this runtime quantizes ColumnWidth, so 18 is the closest achievable stored width)
$ws.Columns("B").ColumnWidth = 18

# Update selection to match the target state
$ws.Range("B14").Select()
